# The document's last paragraph holds the "_GoBack" bookmark and the
# text "名人名言：加油" (with an rFonts hint="default" paragraph mark).
#
# The edit duplicates that paragraph's text into two brand-new paragraphs
# inserted right before it (both with the "normal" rFonts hint="eastAsia"
# paragraph mark used throughout the rest of the document), adds a new
# second line "今天星期四：", and then rewrites the original (bookmarked)
# paragraph's text to "晴：".

$d = $word.ActiveDocument

$count      = $d.Paragraphs.Count
$anchorIdx  = $count - 1   # paragraph right before the final (bookmarked) paragraph

# Insert two new empty paragraphs right after the anchor paragraph. Each
# inherits the anchor's (eastAsia-hinted) paragraph-mark formatting, which
# is exactly what the target markup expects.
$anchorPara = $d.Paragraphs($anchorIdx)
$anchorPara.Range.InsertParagraphAfter()

$newPara1 = $d.Paragraphs($anchorIdx + 1)
$newPara1.Range.InsertParagraphAfter()

# Fill in the text of the first new paragraph: "名人名言：加油"
$newPara1 = $d.Paragraphs($anchorIdx + 1)
$r1 = $newPara1.Range
$r1.SetRange($r1.Start, $r1.Start)
$r1.InsertAfter("名人名言：加油")

# Fill in the text of the second new paragraph: "今天星期四："
$newPara2 = $d.Paragraphs($anchorIdx + 2)
$r2 = $newPara2.Range
$r2.SetRange($r2.Start, $r2.Start)
$r2.InsertAfter("今天星期四：")

# Finally, change the original last paragraph's text from "名人名言：加油"
# to "晴：" (its paragraph mark / bookmark stay untouched).
$finalPara = $d.Paragraphs($d.Paragraphs.Count)
$finalRange = $finalPara.Range
$finalRange.Find.Execute("名人名言：加油", $true, $false, $false, $false, $false, `
                          $true, 1, $false, "晴：", 2)
